$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows after the current last row (274), preserving the
# "D column carries style 4" pattern seen throughout the sheet, while
# stripping the incidental style inheritance on column B.
$ws.Rows.Item(275).Insert()
$ws.Rows.Item(276).Insert()
$ws.Rows.Item(277).Insert()
$ws.Rows.Item(278).Insert()
$ws.Rows.Item(279).Insert()

$ws.Range("B275").ClearFormats()
$ws.Range("B276").ClearFormats()
$ws.Range("B277").ClearFormats()
$ws.Range("B278").ClearFormats()
$ws.Range("B279").ClearFormats()

# Populate cells in the exact order the source strings were first authored,
# so the shared-string table gets appended in that same order.
$ws.Range("B275").Value = "Atwell Island Restoration Project Activities 2000-2010"
$ws.Range("H275").Value = "grey literature"
$ws.Range("J275").Value = "rotational grazing"
$ws.Range("C275").Value = "seeding, micro-topographic relief, fire, grazing"
$ws.Range("L275").Value = "fallowing, disking, micro-topography"
$ws.Range("B276").Value = "Central Valley Project Improvement Act Land Retirement Demonstration Project"
$ws.Range("I276").Value = "retired croplands"
$ws.Range("L276").Value = "salt-remediation products, seed coatings, and mycorrhizal inoculation"
$ws.Range("C276").Value = "irrigation, planting, weed control"
$ws.Range("J276").Value = "future strategy"
$ws.Range("B277").Value = "Efficacy of habitat reclamation for endangered species at the Elk hills oil field, California"
$ws.Range("K277").Value = "habitat for endangered species"
$ws.Range("I277").Value = "oil and gas fields"
$ws.Range("L277").Value = "ripped, disked"
$ws.Range("M277").Value = "undisturbed sites"
$ws.Range("C277").Value = "seeding, ripped, disked, fertilization"
$ws.Range("H277").Value = "grey literature, success-costs"
$ws.Range("B278").Value = "Restoration efforts of Atriplex canescens"
$ws.Range("C278").Value = "grazing, seeding"
$ws.Range("K278").Value = "native species"
$ws.Range("I278").Value = "invasive species"
$ws.Range("K279").Value = "native shrubland habitat"
$ws.Range("M279").Value = "valley saltbush scrub community"
$ws.Range("C279").Value = "seeding, tillage, weed suppression"
$ws.Range("L279").Value = "tillage"
$ws.Range("I279").Value = "solar power plant installation"
$ws.Range("B279").Value = "California Valley Solar Ranch San Luis Obispo County, CA"

# Remaining cells (reuse of existing shared strings / numeric IDs)
# Row 275
$ws.Range("A275").Value = 274
$ws.Range("D275").Value = "field data, App."
$ws.Range("E275").Value = "California, USA"
$ws.Range("F275").Value = "no"
$ws.Range("G275").Value = "NA"
$ws.Range("I275").Value = "agriculture"
$ws.Range("K275").Value = "native vegetation"
$ws.Range("M275").Value = "native vegetation"
$ws.Range("N275").Value = "active restoration"
# Row 276
$ws.Range("A276").Value = 275
$ws.Range("D276").Value = "field data, App."
$ws.Range("E276").Value = "California, USA"
$ws.Range("F276").Value = "no"
$ws.Range("G276").Value = "NA"
$ws.Range("H276").Value = "grey literature"
$ws.Range("K276").Value = "native vegetation"
$ws.Range("M276").Value = "native vegetation"
$ws.Range("N276").Value = "active restoration"
# Row 277
$ws.Range("A277").Value = 276
$ws.Range("D277").Value = "field data"
$ws.Range("E277").Value = "California, USA"
$ws.Range("F277").Value = "no"
$ws.Range("G277").Value = "NA"
$ws.Range("J277").Value = "no"
$ws.Range("N277").Value = "active restoration"
# Row 278
$ws.Range("A278").Value = 277
$ws.Range("D278").Value = "field data"
$ws.Range("E278").Value = "California, USA"
$ws.Range("F278").Value = "no"
$ws.Range("G278").Value = "NA"
$ws.Range("H278").Value = "grey literature"
$ws.Range("J278").Value = "sheep"
$ws.Range("L278").Value = "no"
$ws.Range("M278").Value = "not clear"
$ws.Range("N278").Value = "active restoration"
# Row 279
$ws.Range("A279").Value = 278
$ws.Range("D279").Value = "field data"
$ws.Range("E279").Value = "California, USA"
$ws.Range("F279").Value = "no"
$ws.Range("G279").Value = "NA"
$ws.Range("H279").Value = "grey literature"
$ws.Range("J279").Value = "fenced to livestock"
$ws.Range("N279").Value = "active restoration"
